$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Enemy" placeholder name in column B (rows 9-13, the generic pawn
# entries) is being cleared out so real names can be filled in later -
# replace it with a blank (two-space) placeholder, matching the commit
# "Add Name for pawn".
$ws.Range("B9:B13").Value = "  "

# Row 13 previously carried a highlighted (yellow fill) style that the
# other pawn rows (9-12) don't have; clear it so B13 matches B9:B12.
$ws.Range("B13").Interior.ColorIndex = -4142
$ws.Range("B13").Interior.Pattern = -4142

# Reflect the user's subsequent selection/scroll: they selected the
# newly-cleared name column for the pawn rows and scrolled the view back
# to the left edge (column A) instead of being parked at column U.
$ws.Range("B9:B13").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
